$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 13
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = 7

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 8
